# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme name="Office")
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme name="Red Violet")
# theme2.xml is the theme bound to the (single) Slide Master / the
# presentation's main design, so it's the one the PowerPoint object
# model's Design/Theme/ThemeColorScheme surface exposes.
#
# The authored edit swaps the two themes' content so the deck's visible
# design switches from the pink/purple "Integral" palette to the plain
# "Office" palette. We reproduce that by pushing the "Office Theme"
# color values into the design's ThemeColorScheme, in the fixed slot
# order PowerPoint uses for a color scheme:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
